$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values A1:C1 need to be reordered:
#   before: kitchens_1, bedrooms_1, living_rooms_1
#   after:  living_rooms_1, kitchens_1, bedrooms_1
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "kitchens_1"
$ws.Range("C1").Value = "bedrooms_1"
